{"js": "// Correction bug / slider pagination points \u2014 apply the two text fixes\n// described by the commit \"correction des point de pagination slider\".\n\nconst body = context.document.body;\n\n// --- Edit 1 -----------------------------------------------------------\n// Remove the stray \"components  > \" path segment from the debugging note\n// (\"...du fichier src  > components  > containers  > slider  >  index.js\"\n//  becomes \"...du fichier src  > containers  > slider  >  index.js\").\nconst pathHits = body.search(\"components  > \", { matchCase: false });\npathHits.load(\"items\");\nawait context.sync();\n\nif (pathHits.items.length > 0) {\n  pathHits.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2 -----------------------------------------------------------\n// Fix the wrong \"a\" (should be \"\u00e0\") in:\n// \"Le soucis \u00e9tait que le .map qui servais a g\u00e9n\u00e9r\u00e9 les points du slider...\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Le soucis \u00e9tait que le .map qui servais\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // Narrow the search to this paragraph only: there is also a later,\n  // correctly-spelled \"a\" in \"donc a chaque fois\" that must stay untouched.\n  const aHits = targetParagraph.search(\" a \", { matchCase: true });\n  aHits.load(\"items\");\n  await context.sync();\n\n  if (aHits.items.length > 0) {\n    // Isolate just the \"a\" character (trim the surrounding spaces).\n    const aParts = aHits.items[0].split([\" \"], false, true, true);\n    aParts.load(\"items/text\");\n    await context.sync();\n\n    const aRange = aParts.items[0];\n    // Nudging a character property forces the host to materialise this\n    // sub-range as its own run (matching the 3-run split Word produces\n    // when you retype a single word), instead of silently rewriting the\n    // text inside the existing, larger run.\n    aRange.font.bold = true;\n    await context.sync();\n    aRange.font.bold = false;\n    await context.sync();\n\n    aRange.insertText(\"\u00e0\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Correction bug / slider pagination points -- apply the two text fixes\n# described by the commit \"correction des point de pagination slider\".\n\n$d = $word.ActiveDocument\n\n# --- Edit 1 -------------------------------------------------------------\n# Remove the stray \"components  > \" path segment from the debugging note\n# (\"...du fichier src  > components  > containers  > slider  >  index.js\"\n#  becomes \"...du fichier src  > containers  > slider  >  index.js\").\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"components  > \"\n$find1.Replacement.Text = \"\"\n$find1.Execute([ref]$find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find1.Replacement.Text, 2) | Out-Null\n\n# --- Edit 2 ---------------------------------------------------------------\n# Fix the wrong \"a\" (should be \"\u00e0\") in:\n# \"Le soucis \u00e9tait que le .map qui servais a g\u00e9n\u00e9r\u00e9 les points du slider...\"\n# Only the first \" a \" in that paragraph is wrong -- the later\n# \"donc a chaque fois\" must stay untouched, so locate the word by offset\n# instead of a document-wide Find/Replace.\n$full = $d.Content.Text\n$needle = \"servais a g\u00e9n\u00e9r\u00e9\"\n$needleIdx = $full.IndexOf($needle)\nif ($needleIdx -ge 0) {\n    $aIdx = $needleIdx + 8   # offset of the standalone \"a\" inside the needle\n    $aRange = $d.Range($aIdx, $aIdx + 1)\n    if ($aRange.Text -eq \"a\") {\n        # Replace the character, then nudge a character property: this\n        # forces Word to split this single character off into its own run\n        # (the same 3-run shape you get when retyping just that word),\n        # instead of silently leaving the text inside one merged run.\n        $aRange.Text = \"\u00e0\"\n        $aRange.Font.Bold = 1\n        $aRange.Font.Bold = 0\n    }\n}\n"}
